$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "91.114.52"
$ws.Range("E2").Value = "  +4.76%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.101.29"
$ws.Range("E3").Value = "  +1.78%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.16"
$ws.Range("E5").Value = "  +5.61%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "619.08"

# Row 7
$ws.Range("E7").Value = "  +6.83%  "

# Row 8
$ws.Range("E8").Value = "  +16.84%  "

# Row 9
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.097.87"
$ws.Range("E10").Value = "  +1.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.698"
$ws.Range("E11").Value = "  +22.11%  "

# Row 12
$ws.Range("E12").Value = "  +7.40%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  +10.99%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.40"
$ws.Range("E14").Value = "  +4.07%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.065.85"
$ws.Range("E15").Value = "  +4.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.10"
$ws.Range("E16").Value = "  +7.58%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.673.48"
$ws.Range("E17").Value = "  +1.64%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.082.56"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19
$ws.Range("E19").Value = "  +9.63%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000236"
$ws.Range("E20").Value = "  +15.43%  "

# Row 21
$ws.Range("E21").Value = "  +7.67%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.21"
$ws.Range("E22").Value = "  +4.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.56"
$ws.Range("E23").Value = "  +5.86%  "

# Row 24
$ws.Range("E24").Value = "  +9.38%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.59"
$ws.Range("E25").Value = "  +6.33%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.91"
$ws.Range("E26").Value = "  +7.99%  "

# Row 27
$ws.Range("E27").Value = "  +2.83%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.263.12"
$ws.Range("E28").Value = "  +1.23%  "

# Row 29
$ws.Range("E29").Value = "  +0.15%  "

# Row 30
$ws.Range("E30").Value = "  +13.59%  "

# Row 31
$ws.Range("E31").Value = "  +0.12%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.69"
$ws.Range("E32").Value = "  +10.27%  "

# Row 33
$ws.Range("E33").Value = "  +10.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "518.36"
$ws.Range("E34").Value = "  +6.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.95"
$ws.Range("E35").Value = "  +6.71%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.141"
$ws.Range("E36").Value = "  +1.17%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("E37").Value = "  +5.36%  "

# Row 38
$ws.Range("E38").Value = "  +4.24%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.01"
$ws.Range("E39").Value = "  +6.52%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.30"
$ws.Range("E40").Value = "  +0.87%  "

# Row 41
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.145"
$ws.Range("E42").Value = "  +12.87%  "

# Row 43
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("E44").Value = "  +4.69%  "

# Row 45
$ws.Range("E45").Value = "  +5.60%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0717"
$ws.Range("E46").Value = "  +13.43%  "

# Row 47
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.77"
$ws.Range("E47").Value = "  +0.88%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "142.31"
$ws.Range("E48").Value = "  -2.50%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.27"
$ws.Range("E49").Value = "  +10.31%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000262"
$ws.Range("E50").Value = "  +22.19%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.22"
$ws.Range("E51").Value = "  +10.66%  "
